# Split the "Variables" sheet's 'position' column into two columns:
# 'pivot' (figures/heading/stub) and 'order' (blank/1/1).
# Also clears the old 'type' data column (FIGURES/2MD/2MD) which no
# longer applies once the pivot role moved to its own column.
# Close #124

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new column before A; this shifts old A..G to B..H and keeps
# all their values/styles/widths intact.
$ws.Range("A1").EntireColumn.Insert()

# New column A: "pivot" - replaces the old "position" semantics with the
# figures/heading/stub roles.
$ws.Range("A1").Value = "pivot"
$ws.Range("A2").Value = "figures"
$ws.Range("A3").Value = "heading"
$ws.Range("A4").Value = "stub"

# Column B (previously held "position" data, now relabeled "order")
$ws.Range("B1").Value = "order"
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1

# Column E (shifted from the old D "type" column) keeps its header but
# its FIGURES/2MD/2MD values are no longer used, so clear them.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()

# Update the active sheet/selection to match the new view state.
$ws.Range("D4").Select()
$ws.Activate()

$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("I15").Select()

$ws.Activate()
